$d = $word.ActiveDocument

$d.Content.Find.Execute("39×54=", $true, $false, $false, $false, $false, $true, 1, $false, "46×60=", 2) | Out-Null
$d.Content.Find.Execute("43×15=", $true, $false, $false, $false, $false, $true, 1, $false, "70×93=", 2) | Out-Null
$d.Content.Find.Execute("36×85=", $true, $false, $false, $false, $false, $true, 1, $false, "43×80=", 2) | Out-Null
$d.Content.Find.Execute("68×53=", $true, $false, $false, $false, $false, $true, 1, $false, "88×90=", 2) | Out-Null
$d.Content.Find.Execute("42×89=", $true, $false, $false, $false, $false, $true, 1, $false, "89×66=", 2) | Out-Null
$d.Content.Find.Execute("80×61=", $true, $false, $false, $false, $false, $true, 1, $false, "94×40=", 2) | Out-Null
$d.Content.Find.Execute("96×65=", $true, $false, $false, $false, $false, $true, 1, $false, "67×99=", 2) | Out-Null
$d.Content.Find.Execute("53×72=", $true, $false, $false, $false, $false, $true, 1, $false, "32×21=", 2) | Out-Null
$d.Content.Find.Execute("95×75=", $true, $false, $false, $false, $false, $true, 1, $false, "13×14=", 2) | Out-Null
$d.Content.Find.Execute("52×52=", $true, $false, $false, $false, $false, $true, 1, $false, "99×47=", 2) | Out-Null
$d.Content.Find.Execute("91×82=", $true, $false, $false, $false, $false, $true, 1, $false, "81×99=", 2) | Out-Null
$d.Content.Find.Execute("74×25=", $true, $false, $false, $false, $false, $true, 1, $false, "83×11=", 2) | Out-Null
$d.Content.Find.Execute("21×87=", $true, $false, $false, $false, $false, $true, 1, $false, "33×49=", 2) | Out-Null
$d.Content.Find.Execute("23×59=", $true, $false, $false, $false, $false, $true, 1, $false, "79×19=", 2) | Out-Null
$d.Content.Find.Execute("86×71=", $true, $false, $false, $false, $false, $true, 1, $false, "57×25=", 2) | Out-Null
$d.Content.Find.Execute("95×28=", $true, $false, $false, $false, $false, $true, 1, $false, "80×92=", 2) | Out-Null
$d.Content.Find.Execute("97×40=", $true, $false, $false, $false, $false, $true, 1, $false, "94×67=", 2) | Out-Null
$d.Content.Find.Execute("15×58=", $true, $false, $false, $false, $false, $true, 1, $false, "72×33=", 2) | Out-Null
$d.Content.Find.Execute("78×83=", $true, $false, $false, $false, $false, $true, 1, $false, "91×98=", 2) | Out-Null
$d.Content.Find.Execute("28×12=", $true, $false, $false, $false, $false, $true, 1, $false, "69×69=", 2) | Out-Null
$d.Content.Find.Execute("75×17=", $true, $false, $false, $false, $false, $true, 1, $false, "20×61=", 2) | Out-Null
$d.Content.Find.Execute("74×70=", $true, $false, $false, $false, $false, $true, 1, $false, "26×61=", 2) | Out-Null
$d.Content.Find.Execute("87×65=", $true, $false, $false, $false, $false, $true, 1, $false, "51×76=", 2) | Out-Null
$d.Content.Find.Execute("51×64=", $true, $false, $false, $false, $false, $true, 1, $false, "73×26=", 2) | Out-Null
$d.Content.Find.Execute("34×21=", $true, $false, $false, $false, $false, $true, 1, $false, "42×51=", 2) | Out-Null

Write-Output "done"
